$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.036.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.64%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.466.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.08%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'559.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.77%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'162.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.85%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -1.30%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.465.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.08%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.97%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.61%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  -3.55%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -1.22%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.923.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.88%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'69.042.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.65%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -3.02%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'23.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.09%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.463.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.92%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'10.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.49%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'341.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.03%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'7.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.95%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -2.51%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +0.81%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +0.04%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'67.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.12%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -2.54%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.595.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.05%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.06%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'8.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.59%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0₃0819"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.94%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'7.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.86%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'439.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.25%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +0.09%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -4.07%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -5.35%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'156.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.60%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -0.08%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -0.01%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.92%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'17.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.31%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -3.78%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -2.75%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.94%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'ImmutableX"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.72%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'Stacks"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.84%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.78%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'133.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.05%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'3.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.20%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0720"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.48%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -4.22%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.562"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.98%  "
$ws.Range("E51").Style = "Normal"
